# Updates cryptos list price (D) and 1h volume change (E) columns, and
# swaps the TRON/Polkadot and OKB/Arweave row pairs (rank reordering),
# per the "Updated cryptos list" data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "70.380.91"
$ws.Range("E2").Value = "  +0.99%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.784.13"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.59"
$ws.Range("E5").Value = "  +0.79%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.85"
$ws.Range("E6").Value = "  +2.48%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.789.38"
$ws.Range("E7").Value = "  +1.37%  "

# Row 8: USDC
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.29%  "

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -1.35%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +0.29%  "

# Row 11: Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").Value = "  +1.50%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("E12").Value = "  -3.21%  "

# Row 13: Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.37"
$ws.Range("E13").Value = "  -1.73%  "

# Row 14: ShibaInu
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000259"
$ws.Range("E14").Value = "  +0.80%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.380.55"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.752.46"
$ws.Range("E16").Value = "  -0.19%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "70.307.32"
$ws.Range("E17").Value = "  +0.63%  "

# Row 18: TRON
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.62"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19: Polkadot
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.121"
$ws.Range("E19").Value = "  -2.26%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.70"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "508.30"
$ws.Range("E21").Value = "  -1.98%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.32"
$ws.Range("E22").Value = "  -0.99%  "

# Row 23: Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.731"
$ws.Range("E23").Value = "  -1.28%  "

# Row 24: Fetch.AI
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.64"
$ws.Range("E24").Value = "  +5.68%  "

# Row 25: Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.93"
$ws.Range("E25").Value = "  -2.18%  "

# Row 26: InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.13"
$ws.Range("E26").Value = "  -3.24%  "

# Row 27: RenderToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.38"
$ws.Range("E27").Value = "  +4.21%  "

# Row 28: PEPE
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000136"
$ws.Range("E28").Value = "  +6.69%  "

# Row 29: Dai
$ws.Range("E29").Value = "  +0.35%  "

# Row 30: ImmutableX
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("E30").Value = "  +0.11%  "

# Row 31: PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.95"
$ws.Range("E31").Value = "  +2.84%  "

# Row 32: NEARProtocol
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  +2.23%  "

# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.97"
$ws.Range("E33").Value = "  -2.19%  "

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  -0.69%  "

# Row 35: FirstDigitalUSD
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.29%  "

# Row 36: Mantle
$ws.Range("E36").Value = "  +2.02%  "

# Row 37: Filecoin
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.22"
$ws.Range("E37").Value = "  +0.10%  "

# Row 38: TheGraph
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.352"
$ws.Range("E38").Value = "  +2.46%  "

# Row 39: Kaspa
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").Value = "  +6.24%  "

# Row 40: dogwifhat
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").Value = "  +13.58%  "

# Row 41: Stacks
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.09"
$ws.Range("E41").Value = "  -4.63%  "

# Row 42: OKB
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.89"
$ws.Range("E42").Value = "  +2.78%  "

# Row 43: Arweave
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "50.01"
$ws.Range("E43").Value = "  -2.85%  "

# Row 44: Bittensor
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "440.00"
$ws.Range("E44").Value = "  +2.99%  "

# Row 45: Cosmos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.69"
$ws.Range("E45").Value = "  -1.96%  "

# Row 46: Maker
$ws.Range("D46").Value = "2.984.59"
$ws.Range("E46").Value = "  -3.26%  "

# Row 47: VeChain
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0367"
$ws.Range("E47").Value = "  +0.67%  "

# Row 48: InjectiveProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.63"
$ws.Range("E48").Value = "  -1.15%  "

# Row 49: Monero
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.26"
$ws.Range("E49").Value = "  +2.07%  "

# Row 50: USDe
$ws.Range("E50").Value = "  -0.04%  "

# Row 51: ThetaToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.50"
$ws.Range("E51").Value = "  -0.12%  "
